$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("B1").Value = 16
$ws.Range("C1").Value = 20
$ws.Range("D1").Value = 16
$ws.Range("E1").Value = 20

$ws.Range("B2").Value = 3.3673641978477207
$ws.Range("C2").Value = 6.9552452518524666
$ws.Range("D2").Value = 10.152546027031246
$ws.Range("E2").Value = 8.931617009408086

$ws.Range("B3").Value = 4.6464066055564892
$ws.Range("C3").Value = 5.2343485057408117
$ws.Range("D3").Value = 5.0178586006840877
$ws.Range("E3").Value = 10.127918245401867

[void]$ws.Range("B1:E3").Select()
